$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 196, shifting existing rows 196:300 down to 197:301.
$ws.Rows.Item(196).Insert()

# Fill the new row 196 with the new price-listing record.
# Columns A-C, E-J, L, R repeat the same fixed values used throughout this
# worksheet's block; D, K, M, N, O, P, Q, S, T carry the new record's data.
$ws.Range("A196").Value = 7
$ws.Range("B196").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C196").Value = "Ñuble"
$ws.Range("D196").Value = 44518
$ws.Range("E196").Value = 16
$ws.Range("F196").Value = "Fruta"
$ws.Range("G196").Value = 100102
$ws.Range("H196").Value = "Cítricos"
$ws.Range("I196").Value = 100102005
$ws.Range("J196").Value = "Naranja"
$ws.Range("K196").Value = "Lane Late"
$ws.Range("L196").Value = "Primera"
$ws.Range("M196").Value = 160
$ws.Range("N196").Value = 8500
$ws.Range("O196").Value = 9000
$ws.Range("P196").Value = 8750
$ws.Range("Q196").Value = "$/bandeja 15 kilos granel"
$ws.Range("R196").Value = "Región de O'Higgins"
$ws.Range("S196").Value = 583
$ws.Range("T196").Value = 15
